$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 29 with the latest ranking update timestamp and placeholder "-" values
$ws.Range("A29").Value = "2025/12/03 21:00"
$ws.Range("B29").Value = "-"
$ws.Range("C29").Value = "-"
$ws.Range("D29").Value = "-"
$ws.Range("E29").Value = "-"
$ws.Range("F29").Value = "-"
$ws.Range("G29").Value = "-"
